# Update cryptos list - apply scraped price/volume changes and
# the Stellar/ImmutableX row swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set D (Price) and E (Volume 1h) columns for a row, keeping
# values as plain text (the Price/Volume columns are textual, and many
# values look numeric so Excel would otherwise coerce them).
function Set-Row($row, $price, $volume) {
    if ($price -ne $null) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $price
    }
    if ($volume -ne $null) {
        $ecell = $ws.Cells.Item($row, 5)
        $ecell.NumberFormat = "@"
        $ecell.Value = $volume
    }
}

Set-Row 2  "28.206.02"   "  +0.58%  "
Set-Row 3  "1.872.32"    "  +4.00%  "
Set-Row 4  $null          "  -0.04%  "
Set-Row 5  "311.53"      "  +0.44%  "
Set-Row 6  $null          "  -0.09%  "
Set-Row 7  "0.4998"      "  -1.52%  "
Set-Row 8  "0.3910"      "  +2.05%  "
Set-Row 9  "0.09540"     "  +23.60%  "
Set-Row 10 "1.140"       "  +4.25%  "
Set-Row 11 "41.04"       "  +0.79%  "
Set-Row 12 "6.472"       "  +1.81%  "
Set-Row 13 "20.98"       "  +3.20%  "
Set-Row 14 "1.877.82"    "  +4.07%  "
Set-Row 15 $null          "  -0.04%  "
Set-Row 16 "7.391"       "  +1.72%  "
Set-Row 17 $null          "  +4.88%  "
Set-Row 18 "93.30"       "  +1.39%  "
Set-Row 19 "0.06619"     "  +0.78%  "
Set-Row 20 "17.50"       "  +1.55%  "
Set-Row 21 $null          "  -0.04%  "
Set-Row 22 "6.145"       "  +2.85%  "
Set-Row 23 "28.265.45"   "  +0.76%  "
Set-Row 24 "11.33"       "  +2.76%  "
Set-Row 25 $null          "  +2.49%  "
Set-Row 26 "2.552"       "  +5.61%  "
Set-Row 27 "2.090.05"    "  +3.92%  "
Set-Row 28 "21.17"       "  +4.74%  "
Set-Row 29 "157.24"      "  -1.32%  "
Set-Row 30 "127.81"      "  +0.53%  "

# Row 31 / Row 32 swap: Stellar <-> ImmutableX
$ws.Cells.Item(31, 2).Value = "ImmutableX"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c31d = $ws.Cells.Item(31, 4)
$c31d.NumberFormat = "@"
$c31d.Value = "1.066"
$c31e = $ws.Cells.Item(31, 5)
$c31e.NumberFormat = "@"
$c31e.Value = "  +2.03%  "

$ws.Cells.Item(32, 2).Value = "Stellar"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c32d = $ws.Cells.Item(32, 4)
$c32d.NumberFormat = "@"
$c32d.Value = "0.1054"
$c32e = $ws.Cells.Item(32, 5)
$c32e.NumberFormat = "@"
$c32e.Value = "  -3.37%  "

Set-Row 33 "5.631"       "  +1.74%  "
Set-Row 34 "3.626"       "  -0.53%  "
Set-Row 35 "0.06760"     "  -2.65%  "
Set-Row 36 "9.567"       "  +5.14%  "
Set-Row 37 "0.02395"     "  +2.91%  "
Set-Row 38 "0.2175"      "  +0.48%  "
Set-Row 39 "11.50"       "  +0.58%  "
Set-Row 40 "4.977"       "  -0.30%  "
Set-Row 41 "0.6305"      "  +3.39%  "
Set-Row 42 $null          "  +2.35%  "
Set-Row 43 "1.001"       "  -0.01%  "
Set-Row 44 "13.56"       "  +2.48%  "
Set-Row 45 "0.6038"      "  +2.73%  "
Set-Row 46 "3.659"       "  -1.13%  "
Set-Row 47 "1.260"       "  -2.36%  "
Set-Row 48 "123.94"      "  -1.42%  "
Set-Row 49 "1.984"       "  +2.93%  "
Set-Row 50 "1.194"       "  +0.95%  "
Set-Row 51 "0.06845"     "  +1.89%  "
